$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HMI")

$ws.Range("D8").Value = 304100
$ws.Range("E8").Value = 231000
$ws.Range("F8").Value = 133000

$ws.Range("D9").Value = 230700
$ws.Range("E9").Value = 190000
$ws.Range("F9").Value = 116600

$ws.Range("D10").Value = 73400
$ws.Range("E10").Value = 41000
$ws.Range("F10").Value = 16400

$ws.Range("D12").Value = 22800
$ws.Range("E12").Value = 19600
$ws.Range("F12").Value = 9100

$ws.Range("D17").Value = 277100
$ws.Range("E17").Value = 229000
$ws.Range("F17").Value = 139000

$ws.Range("D18").Value = 27000
$ws.Range("E18").Value = 2000
$ws.Range("F18").Value = -6000

$ws.Range("D20").Value = 1500

$ws.Range("D21").Value = 28800
$ws.Range("E21").Value = 4300
$ws.Range("F21").Value = "NA"

$ws.Range("D23").Value = 28500
$ws.Range("E23").Value = 4300
$ws.Range("F23").Value = -5800

$ws.Range("D24").Value = 3500
$ws.Range("E24").Value = 500

$ws.Range("D26").Value = 25000
$ws.Range("E26").Value = 3800
$ws.Range("F26").Value = -5600

$ws.Range("D27").Value = 7500
$ws.Range("F27").Value = -9100

$ws.Range("D32").Value = -1500

$ws.Range("D33").Value = 6800
$ws.Range("F33").Value = -9100

$ws.Range("D35").Value = 6800
$ws.Range("F35").Value = -9100

$ws.Range("D41").Value = 54400
$ws.Range("E41").Value = 22700
$ws.Range("F41").Value = 32600

$ws.Range("E42").Value = 1400

$ws.Range("D43").Value = 95100
$ws.Range("E43").Value = 74000
$ws.Range("F43").Value = 30500

$ws.Range("D44").Value = 37100
$ws.Range("E44").Value = 28500
$ws.Range("F44").Value = 13300

$ws.Range("D46").Value = 192200
$ws.Range("E46").Value = 127600
$ws.Range("F46").Value = 77100

$ws.Range("D47").Value = 12700
$ws.Range("E47").Value = 11600

$ws.Range("D48").Value = 4300

$ws.Range("D49").Value = 1700

$ws.Range("D52").Value = 6700
$ws.Range("E52").Value = 3400

$ws.Range("D54").Value = 217500
$ws.Range("E54").Value = 144400
$ws.Range("F54").Value = 78500

$ws.Range("D57").Value = 106300
$ws.Range("E57").Value = 81300
$ws.Range("F57").Value = 38600

$ws.Range("D58").Value = 5200
$ws.Range("E58").Value = 1900

$ws.Range("D59").Value = 18700
$ws.Range("E59").Value = 11000

$ws.Range("D60").Value = 130200
$ws.Range("E60").Value = 94100
$ws.Range("F60").Value = 41200

$ws.Range("D62").Value = 1600

$ws.Range("D66").Value = 132100
$ws.Range("E66").Value = 94100
$ws.Range("F66").Value = 41200

$ws.Range("D70").Value = 51900
$ws.Range("E70").Value = 45800
$ws.Range("F70").Value = 40400

$ws.Range("D72").Value = 19500
$ws.Range("E72").Value = -5400
$ws.Range("F72").Value = -9000

$ws.Range("D76").Value = 33500

$ws.Range("D81").Value = 6800
$ws.Range("F81").Value = -9100

$ws.Range("D83").Value = 400
$ws.Range("E83").Value = 0
$ws.Range("F83").Value = "NA"

$ws.Range("D89").Value = 2600
$ws.Range("E89").Value = -1000
$ws.Range("F89").Value = "NA"

$ws.Range("D91").Value = -1500
$ws.Range("E91").Value = -400
$ws.Range("F91").Value = "NA"

$ws.Range("D94").Value = -14800
$ws.Range("E94").Value = -700
$ws.Range("F94").Value = "NA"

$ws.Range("D100").Value = 1500
$ws.Range("E100").Value = 31800
$ws.Range("F100").Value = "NA"

$ws.Range("D101").Value = 800
$ws.Range("E101").Value = 1500
$ws.Range("F101").Value = "NA"

$ws.Range("D102").Value = -9900
$ws.Range("E102").Value = 31600
$ws.Range("F102").Value = "NA"
